# Applies the DaCapo Shenandoah GC heap-1G stats/docx fix.
#
# The first table column had several single-value rows (rows 1-4 and
# 6-12, 1-based) whose numbers were stale, while three later rows
# (44-46) each held a full tab-separated summary line (10 values) that
# really only needed to show the single overall total. This edit:
#   - rewrites the stale single values in rows 1-4 and 6-12
#   - collapses the tab-separated rows 44-46 down to a single value
#     (re-using the original values that used to sit in rows 1-3)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $expectedOld, $newText) {
    $cell = $table.Cell($rowIndex, 1)
    $r = $cell.Range
    # Trim trailing end-of-cell marker so we don't clobber the cell/row end.
    $r.End = $r.End - 1
    if ($r.Text -ne $expectedOld) {
        throw "Row $rowIndex unexpected text: '$($r.Text)' (expected '$expectedOld')"
    }
    $r.Text = $newText
}

Set-CellText $t 1 "99.98" "0M"
Set-CellText $t 2 "0.08" "0M"
Set-CellText $t 3 "412" "0M"
Set-CellText $t 4 "196" "418"

Set-CellText $t 6 "0.00034" "0.00085"
Set-CellText $t 7 "0.00012" "0.00016"
Set-CellText $t 8 "0.00004" "0.00005"
Set-CellText $t 9 "0.00014" "0.00027"
Set-CellText $t 10 "0.00016" "0.00034"
Set-CellText $t 11 "0.00019" "0.00045"
Set-CellText $t 12 "0.02474" "0.08023"

Set-CellText $t 44 "111`t0.00018`t0.00085`t0.00037`t0.00014`t0.00027`t0.00034`t0.00045`t0.04159`t100.0" "99.98"
Set-CellText $t 45 "26`t0.00003`t0.00009`t0.00005`t0.00001`t0.00004`t0.00005`t0.00006`t0.00133`t100.0" "0.08"
Set-CellText $t 46 "85`t0.00009`t0.00026`t0.00015`t0.00003`t0.00012`t0.00014`t0.00016`t0.01257`t100.0" "412"
